$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "What happened in project:" / "What happened for course:" notes (column C) ---
# Each of these appends a new line to the existing text.

$ws.Range("C4").Value = "Finished Proposal up to standard that I am happy" + [char]10 + "Itial understanding of the process involved throughout the project and what extra learning needs to be done. (c# web services)" + [char]10 + "Will have to relearn angularJS and how to utilise that"

$ws.Range("C6").Value = "Finish Risk analysis and sent through my final proposal to Mike" + [char]10 + "Left with lots more questions but made big advancements"

$ws.Range("C7").Value = "Start detailed planning and initial diagrams" + [char]10 + "Connect through to Dev server" + [char]10 + "Test connections of Tourplan Webservice as well as to GIATA" + [char]10 + "Understand the connection between angularjs and c#"

$ws.Range("C8").Value = "Connecting through to the Development server on Tourplany system" + [char]10 + "The initial understanding of angularJS"

# C5 text is unchanged: "Finished Proposal up to standard that I am happy"

# --- New column D: "Notes:" entries for each day ---

$ws.Range("D4").Value = "Some progress on using angularJS and how to incorperate that into my application" + [char]10 + "Initial Understanding of the data that needs to flow between entities"
$ws.Range("D4").WrapText = $true
$ws.Range("D4").VerticalAlignment = -4160

$ws.Range("D5").Value = "Wrote in diary + Time schedule" + [char]10 + "Had Buisness analysis that meant I had to go into ara "
$ws.Range("D5").WrapText = $true

$ws.Range("D6").Value = "Connected to Dev server" + [char]10 + "Beginning to understand the relationships" + [char]10 + "Not mauch initial planning "
$ws.Range("D6").WrapText = $true

$ws.Range("D7").Value = "More comprehensive planning with my general idea layed out with work componentry to say it is possible" + [char]10 + " Have my proposal checked off"
$ws.Range("D7").WrapText = $true
$ws.Range("D7").VerticalAlignment = -4160

$ws.Range("D8").Value = "Hard getting my head around angular JS again." + [char]10 + "Not a complete understanding of what is required but have enough information to get on with some work. Will meet early next week to discuss my progress"
$ws.Range("D8").WrapText = $true

# --- Column D width ---
$ws.Columns("D").ColumnWidth = 32.71

# --- Row heights to fit the longer text ---
$ws.Rows(4).RowHeight = 138.75
$ws.Rows(6).RowHeight = 60.75
$ws.Rows(7).RowHeight = 105.75
$ws.Rows(8).RowHeight = 105.75

# --- Selection / view state ---
$ws.Range("D17:D19").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
